# Applies the "revise cheatsheets for api changes" edits:
#  - ds.crop / ds.select  -> ds.subset  (Subsetting data box)
#  - ds.drop([...])       -> ds.drop(variables = [...])
#  - ds.merge()           -> ds.merge("variables")
#  - ds.merge_time()      -> ds.merge("time")
#  - ds.set_missing       -> ds.as_missing

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)

# --- Slide 1 : "Subsetting data" box (TextBox 70) ---
$subsetBox = $s1.Shapes.Item(5)
$subsetText = $subsetBox.TextFrame.TextRange

$subsetText.Paragraphs(1).Text = "ds.subset(lon = [lon_min, lon_max],"
$subsetText.Paragraphs(4).Text = "ds.subset(variables = [var1, var2])"
$subsetText.Paragraphs(6).Text = "ds.subset(years = [2000, 2001])" + [char]9
$subsetText.Paragraphs(8).Text = "ds.subset(months = [5, 6])" + [char]9
$subsetText.Paragraphs(10).Text = "ds.drop(variables = [" + [char]8216 + "var1" + [char]8217 + ", " + [char]8216 + "var2])"

# --- Slide 1 : "Merging methods" box (TextBox 103) ---
$mergeBox = $s1.Shapes.Item(17)
$mergeText = $mergeBox.TextFrame.TextRange

$mergeText.Paragraphs(1).Text = "ds.merge(" + [char]8220 + "variables" + [char]8221 + ")"
$mergeText.Paragraphs(3).Text = "ds.merge(" + [char]8220 + "time" + [char]8221 + ")"

# --- Slide 2 : "Random hacks" box (TextBox 24) ---
$hacksBox = $s2.Shapes.Item(6)
$hacksText = $hacksBox.TextFrame.TextRange

$hacksText.Paragraphs(5).Text = "ds.as_missing([0, 100])"
